# Rename PFRCost as REGCV1Cost for clarification
$wb = $excel.ActiveWorkbook

# Rename the "PFRCost" sheet to "REGCV1Cost"
$ws = $wb.Worksheets.Item("PFRCost")
$ws.Name = "REGCV1Cost"

# Make the renamed sheet the active sheet/tab, with the given selection
$ws.Activate()
$ws.Range("J27").Select()
